# Attenuation.xlsx edit
# Commit: "A star destroyer model added"
#
# The "Linear" and "Quadratic" coefficients (M3, N3) are reduced from 0.2
# to 0.05, which changes the attenuation curve so it no longer saturates
# the brightness column at 1 for the first few rows -- the falloff
# becomes gentler. Two helper cells are also added (R3/R4) with a small
# note (S3) flagging that the previous value was "Too big".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core parameter change: Linear (M3) and Quadratic (N3) coefficients ---
$ws.Range("M3").Value = 0.05
$ws.Range("N3").Value = 0.05

# --- New helper notes added off to the side of the model ---
$ws.Range("R3").Value = 0.5
$ws.Range("R4").Value = 0.25
$ws.Range("S3").Value = "Too big"

# Re-enter the "Distance" running-total formula across A4:A27 as one
# block so Excel consolidates it into a single shared formula (matching
# what a resave of the workbook produces) instead of 24 separate ones.
$ws.Range("A4:A27").Formula = '=A3+$K$3'

# Put the selection where the author apparently left it and restore a
# sensible zoom level for the sheet view.
[void]$ws.Range("N4").Select()
$excel.ActiveWindow.Zoom = 100
